$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.726.87"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.601.26"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.74"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.80"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.595.31"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.97"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.689.80"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.51"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.05"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.37"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.292.63"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.598"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("E38").Value = "  +6.02%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.828"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.781"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.99"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.738.57"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.42"
$ws.Range("E51").Value = "  +0.18%  "
